# Updates the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# tracker sheet with freshly scraped values (GitHub Actions refresh).
#
# The Price column holds plain text (prices are dot-grouped, e.g.
# "26.008.18", which isn't a real number) so for any D-cell whose new
# value happens to *look* like a normal decimal number (e.g. "18.47"),
# we force the cell to Text format first - otherwise Excel would
# auto-convert the assignment into a numeric value (and mangle values
# like "0.530" -> 0.53). The number format is reset back to the
# worksheet's default ("Normal" style) right after, so no cell ends up
# visibly reformatted.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.999.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.633.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("E6").Value = "  -0.62%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -2.06%  "
$ws.Range("E9").Value = "  -2.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.94%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.861.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.633.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.530"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.005.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("E17").Value = "  -2.75%  "
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "190.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("E23").Value = "  -2.08%  "
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.31"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("E27").Value = "  -2.06%  "
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("E31").Value = "  -3.40%  "
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("E33").Value = "  -4.08%  "
$ws.Range("E34").Value = "  -1.55%  "
$ws.Range("E35").Value = "  -2.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.134.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("E37").Value = "  -3.74%  "
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("E40").Value = "  -1.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("E43").Value = "  -4.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.772.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("E45").Value = "  -0.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.22"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.70%  "
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("E48").Value = "  +2.16%  "
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.60%  "
$ws.Range("E51").Value = "  +0.08%  "
